$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.932.69"
$ws.Range("E2").Value = "  -0.69%  "

$ws.Range("D3").Value = "1.863.57"
$ws.Range("E3").Value = "  -0.37%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.85"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5070"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.27%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3631"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07173"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8957"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.70"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.00%  "

$ws.Range("D12").Value = "1.857.69"
$ws.Range("E12").Value = "  -0.69%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07445"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.59%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.239"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.67%  "

$ws.Range("E16").Value = "  -0.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008492"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.0000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.08%  "

$ws.Range("D20").Value = "26.973.23"
$ws.Range("E20").Value = "  -0.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.023"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.30%  "

$ws.Range("D22").Value = "2.089.85"
$ws.Range("E22").Value = "  -1.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.435"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.77%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.794"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.70%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.75%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.068"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.670"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.93%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.676"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09235"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.77%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05086"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.992"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7476"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.152"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.88%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.289"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.523"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.04%  "

$ws.Range("E39").Value = "  -1.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.081"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5378"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "117.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.50%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.498"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.91%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.569"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.02%  "

$ws.Range("E45").Value = "  -0.70%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4659"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9995"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.564"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.47%  "
